$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (long name) on both sheets
$wsInput.Range("B1").Value = "4218-RBI-EI-DB-DL-REC-RNI-FEE-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"
$wsOutput.Range("B1").Value = "4218-RBI-EI-DB-DL-REC-RNI-FEE-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"

# Update shortname to text value "421e"
$wsInput.Range("B2").Value = "421e"

# Move selection on input sheet to B1
$wsInput.Range("B1").Select() | Out-Null

# Move selection on output sheet to B1 and make it the active/selected tab
$wsOutput.Range("B1").Select() | Out-Null
$wsOutput.Activate() | Out-Null
